$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1098.9
$ws.Range("I2").Value = 149.6
$ws.Range("K2").Value = 149.6
$ws.Range("M2").Value = -36.59999999999999

# Row 12
$ws.Range("H12").Value = 7271.143
$ws.Range("I12").Value = 3474.75
$ws.Range("K12").Value = 3474.75
$ws.Range("M12").Value = -3304.75

# Row 113
$ws.Range("H113").Value = 3847.4
$ws.Range("I113").Value = 2829.6667
$ws.Range("K113").Value = 2829.6667
$ws.Range("M113").Value = 424.3332999999998

# Row 137
$ws.Range("H137").Value = 3748.8823
$ws.Range("I137").Value = 2356.111
$ws.Range("J137").Value = 5315.75
$ws.Range("K137").Value = 7068.333
$ws.Range("L137").Value = 15947.25
$ws.Range("M137").Value = -4518.333
$ws.Range("N137").Value = -21047.25

# Row 138
$ws.Range("H138").Value = 4495.7886
$ws.Range("I138").Value = 2449.7693
$ws.Range("J138").Value = 6541.8076
$ws.Range("K138").Value = 7349.3079
$ws.Range("L138").Value = 19625.4228
$ws.Range("M138").Value = -2209.3079
$ws.Range("N138").Value = -29905.4228

# Row 141
$ws.Range("H141").Value = 2565.4092
$ws.Range("I141").Value = 1299.7576
$ws.Range("K141").Value = 3899.2728
$ws.Range("M141").Value = 1280.7272

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 175.72728
$ws.Range("I5").Value = 175.72728
$ws.Range("K5").Value = 175.72728
$ws.Range("M5").Value = -63.72728000000001

# Row 61
$ws.Range("H61").Value = 3948.6667
$ws.Range("I61").Value = 2831.1875
$ws.Range("K61").Value = 2831.1875
$ws.Range("M61").Value = -2619.1875

# Row 74
$ws.Range("H74").Value = 2003.4445
$ws.Range("I74").Value = 1966.375
$ws.Range("K74").Value = 1966.375
$ws.Range("M74").Value = -1092.375

# Row 77
$ws.Range("H77").Value = 2003.4445
$ws.Range("I77").Value = 1966.375
$ws.Range("K77").Value = 9831.875
$ws.Range("M77").Value = -5463.875

# Row 110
$ws.Range("H110").Value = 148043.94
$ws.Range("I110").Value = 193213.03
$ws.Range("K110").Value = 193213.03
$ws.Range("M110").Value = -191168.03

# Row 122
$ws.Range("H122").Value = 5431.4375
$ws.Range("I122").Value = 3742.1667
$ws.Range("K122").Value = 11226.5001
$ws.Range("M122").Value = -8776.500100000001

# Row 132
$ws.Range("H132").Value = 3494.9443
$ws.Range("I132").Value = 2122.76
$ws.Range("J132").Value = 6613.5454
$ws.Range("K132").Value = 6368.280000000001
$ws.Range("L132").Value = 19840.6362
$ws.Range("M132").Value = -3838.280000000001
$ws.Range("N132").Value = -24900.6362

# Row 136
$ws.Range("H136").Value = 3948.6667
$ws.Range("I136").Value = 2831.1875
$ws.Range("K136").Value = 8493.5625
$ws.Range("M136").Value = -5943.5625

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 175.72728
$ws.Range("I4").Value = 175.72728
$ws.Range("K4").Value = 175.72728
$ws.Range("M4").Value = -60.72728000000001

# Row 5
$ws.Range("H5").Value = 749.5
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 999
$ws.Range("M5").Value = -387
$ws.Range("N5").Value = -1225

# Row 107
$ws.Range("H107").Value = 2063.138
$ws.Range("I107").Value = 2024.2174
$ws.Range("K107").Value = 2024.2174
$ws.Range("M107").Value = -104.2174

# Row 134
$ws.Range("H134").Value = 29467.44
$ws.Range("I134").Value = 3290.087
$ws.Range("K134").Value = 9870.261
$ws.Range("M134").Value = -7335.261

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 718100.5
$ws.Range("I31").Value = 835700.5600000001
$ws.Range("J31").Value = 12500
$ws.Range("K31").Value = 835700.5600000001
$ws.Range("L31").Value = 12500
$ws.Range("M31").Value = -835405.5600000001
$ws.Range("N31").Value = -13090

# Row 34
$ws.Range("H34").Value = 718100.5
$ws.Range("I34").Value = 835700.5600000001
$ws.Range("J34").Value = 12500
$ws.Range("K34").Value = 835700.5600000001
$ws.Range("L34").Value = 12500
$ws.Range("M34").Value = -835498.5600000001
$ws.Range("N34").Value = -12904

# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

# Row 58
$ws.Range("H58").Value = 248986.1
$ws.Range("I58").Value = 837177.4399999999
$ws.Range("K58").Value = 837177.4399999999
$ws.Range("M58").Value = -836974.4399999999

# Row 105
$ws.Range("H105").Value = 1493.3334
$ws.Range("I105").Value = 822.44446
$ws.Range("K105").Value = 822.44446
$ws.Range("M105").Value = 924.55554

# Row 132
$ws.Range("H132").Value = 2923.6216
$ws.Range("I132").Value = 1996.4286
$ws.Range("K132").Value = 5989.2858
$ws.Range("M132").Value = -3459.2858

# Row 134
$ws.Range("H134").Value = 196222.06
$ws.Range("I134").Value = 2498.3103
$ws.Range("J134").Value = 440482.44
$ws.Range("K134").Value = 7494.9309
$ws.Range("L134").Value = 1321447.32
$ws.Range("M134").Value = -4959.9309
$ws.Range("N134").Value = -1326517.32

# Row 136
$ws.Range("H136").Value = 248986.1
$ws.Range("I136").Value = 837177.4399999999
$ws.Range("K136").Value = 2511532.32
$ws.Range("M136").Value = -2508982.32

$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 1776
$ws.Range("I51").Value = 1718.1666
$ws.Range("J51").Value = 1949.5
$ws.Range("K51").Value = 5154.4998
$ws.Range("L51").Value = 5848.5
$ws.Range("M51").Value = -4694.4998
$ws.Range("N51").Value = -6768.5

# Row 68
$ws.Range("H68").Value = 2503583.5
$ws.Range("I68").Value = 3336888.8
$ws.Range("J68").Value = 2003600.6
$ws.Range("K68").Value = 10010666.4
$ws.Range("L68").Value = 6010801.800000001
$ws.Range("M68").Value = -10009855.4
$ws.Range("N68").Value = -6012423.800000001

# Row 71
$ws.Range("H71").Value = 2503583.5
$ws.Range("I71").Value = 3336888.8
$ws.Range("J71").Value = 2003600.6
$ws.Range("K71").Value = 30031999.2
$ws.Range("L71").Value = 18032405.4
$ws.Range("M71").Value = -30027943.2
$ws.Range("N71").Value = -18040517.4

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 66664.336
$ws.Range("I24").Value = 100000
$ws.Range("K24").Value = 100000
$ws.Range("M24").Value = -99827

# Row 122
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

# Row 132
$ws.Range("H132").Value = 271670.22
$ws.Range("I132").Value = 373227.66
$ws.Range("J132").Value = 75809.42999999999
$ws.Range("K132").Value = 1119682.98
$ws.Range("L132").Value = 227428.29
$ws.Range("M132").Value = -1117152.98
$ws.Range("N132").Value = -232488.29

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1731.125
$ws.Range("I100").Value = 1735
$ws.Range("K100").Value = 1735
$ws.Range("M100").Value = -1194

# Row 132
$ws.Range("H132").Value = 3406.347
$ws.Range("I132").Value = 2881.6511
$ws.Range("K132").Value = 8644.953300000001
$ws.Range("M132").Value = -6114.953300000001

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 7000
$ws.Range("J18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7346

# Row 20
$ws.Range("H20").Value = 40000
$ws.Range("I20").Value = 40000
$ws.Range("K20").Value = 40000
$ws.Range("M20").Value = -39760

# Row 29
$ws.Range("H29").Value = 333353340
$ws.Range("J29").Value = 333353340
$ws.Range("L29").Value = 333353340
$ws.Range("N29").Value = -333353920

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

# Row 75
$ws.Range("H75").Value = 37166.668
$ws.Range("J75").Value = 37166.668
$ws.Range("L75").Value = 37166.668
$ws.Range("N75").Value = -39038.668

# Row 78
$ws.Range("H78").Value = 37166.668
$ws.Range("J78").Value = 37166.668
$ws.Range("L78").Value = 111500.004
$ws.Range("N78").Value = -120860.004

# Row 99
$ws.Range("H99").Value = 45000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 45000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 45000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -50990

# Row 122
$ws.Range("H122").Value = 40004064
$ws.Range("I122").Value = 66669176
$ws.Range("J122").Value = 6400.4
$ws.Range("K122").Value = 200007528
$ws.Range("L122").Value = 19201.2
$ws.Range("M122").Value = -200005078
$ws.Range("N122").Value = -24101.2

# Row 132
$ws.Range("H132").Value = 14445.867
$ws.Range("I132").Value = 1771.5106
$ws.Range("K132").Value = 5314.531800000001
$ws.Range("M132").Value = -2784.531800000001

# Row 136
$ws.Range("H136").Value = 56310.445
$ws.Range("I136").Value = 12773.667
$ws.Range("K136").Value = 38321.001
$ws.Range("M136").Value = -35771.001
